# Updates loading_percent values for rows 2-25 (case with 380 kV) across
# columns B, C, D, E, G, I, K, matching the new simulation results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data["B2"] = 11.19518837836857
$data["C2"] = 13.33274311299466
$data["D2"] = 6.036471631099489
$data["E2"] = 16.17844495477147
$data["G2"] = 3.716816261039574
$data["I2"] = 46.98928765420878
$data["K2"] = 13.77841183086004
$data["B3"] = 11.17945846270876
$data["C3"] = 12.87446737151539
$data["D3"] = 5.928187649198708
$data["E3"] = 15.27763314239106
$data["G3"] = 3.72143590021131
$data["I3"] = 45.43681046224155
$data["K3"] = 13.68413294089438
$data["B4"] = 11.17785992256238
$data["C4"] = 12.5878789107174
$data["D4"] = 5.862738700480048
$data["E4"] = 14.70269768799155
$data["G4"] = 3.724406447189796
$data["I4"] = 44.45500768617126
$data["K4"] = 13.63231296142133
$data["B5"] = 11.17922263987644
$data["C5"] = 12.46996264395679
$data["D5"] = 5.836365576645203
$data["E5"] = 14.46319474594954
$data["G5"] = 3.725650859537178
$data["I5"] = 44.04814025889193
$data["K5"] = 13.6127383608457
$data["B6"] = 11.17957012473237
$data["C6"] = 12.4503200825657
$data["D6"] = 5.832005400147928
$data["E6"] = 14.42311983037235
$data["G6"] = 3.725859545333009
$data["I6"] = 43.98018380810109
$data["K6"] = 13.60958159494585
$data["B7"] = 11.17787016485491
$data["C7"] = 12.58629297416406
$data["D7"] = 5.862381768816507
$data["E7"] = 14.69948836519361
$data["G7"] = 3.724423092313541
$data["I7"] = 44.44954738989045
$data["K7"] = 13.63204270656754
$data["B8"] = 11.18808557734818
$data["C8"] = 13.17591325036929
$data["D8"] = 5.998942104584263
$data["E8"] = 15.87252782685406
$data["G8"] = 3.718381393727388
$data["I8"] = 46.46015884965404
$data["K8"] = 13.74465619185391
$data["B9"] = 11.27242842816386
$data["C9"] = 14.28372079499504
$data["D9"] = 6.273346988098577
$data["E9"] = 17.98982067632706
$data["G9"] = 3.707589286507261
$data["I9"] = 50.16014188053964
$data["K9"] = 14.01273786689276
$data["B10"] = 11.37380102451482
$data["C10"] = 15.05961644868051
$data["D10"] = 6.476830641339252
$data["E10"] = 19.58460887140967
$data["G10"] = 3.700292409467271
$data["I10"] = 52.71159191253102
$data["K10"] = 14.2370826294095
$data["B11"] = 11.42842978139915
$data["C11"] = 15.40284085988572
$data["D11"] = 6.569368021013387
$data["E11"] = 20.28004451057869
$data["G11"] = 3.697107639164666
$data["I11"] = 53.83277449429655
$data["K11"] = 14.34472995972216
$data["B12"] = 11.45033060393105
$data["C12"] = 15.53130507352366
$data["D12"] = 6.604370341291475
$data["E12"] = 20.53764025668208
$data["G12"] = 3.695920811813225
$data["I12"] = 54.25142682418722
$data["K12"] = 14.38626375032252
$data["B13"] = 11.44556007636047
$data["C13"] = 15.50370662637749
$data["D13"] = 6.596834274323292
$data["E13"] = 20.48241754865038
$data["G13"] = 3.6961755663592
$data["I13"] = 54.16152909938873
$data["K13"] = 14.37728501661518
$data["B14"] = 11.43020731522042
$data["C14"] = 15.41344046958874
$data["D14"] = 6.572248673992219
$data["E14"] = 20.30135209712868
$data["G14"] = 3.697009614892395
$data["I14"] = 53.86733707690187
$data["K14"] = 14.34813173740464
$data["B15"] = 11.42096104767676
$data["C15"] = 15.35795055821692
$data["D15"] = 6.557183102457896
$data["E15"] = 20.18969652658837
$data["G15"] = 3.6975229855798
$data["I15"] = 53.68635909931058
$data["K15"] = 14.33037375863331
$data["B16"] = 11.37040144714465
$data["C16"] = 15.03698054671242
$data["D16"] = 6.470779466292374
$data["E16"] = 19.53835368941189
$data["G16"] = 3.700503235015666
$data["I16"] = 52.63750449049469
$data["K16"] = 14.23015728045499
$data["B17"] = 11.34155938989409
$data["C17"] = 14.83750187384691
$data["D17"] = 6.417740641700538
$data["E17"] = 19.12849190017381
$data["G17"] = 3.702365870999141
$data["I17"] = 51.98377474564732
$data["K17"] = 14.17008632625375
$data["B18"] = 11.32577243395015
$data["C18"] = 14.72185634086774
$data["D18"] = 6.387233478731978
$data["E18"] = 18.88895842463062
$data["G18"] = 3.703449891539042
$data["I18"] = 51.60406239483179
$data["K18"] = 14.1360635957903
$data["B19"] = 11.32056527117537
$data["C19"] = 14.68254783508937
$data["D19"] = 6.376905274916819
$data["E19"] = 18.80720407624201
$data["G19"] = 3.703819105851589
$data["I19"] = 51.47486982370724
$data["K19"] = 14.12463585760592
$data["B20"] = 11.3445467093398
$data["C20"] = 14.85883179159204
$data["D20"] = 6.423387042708057
$data["E20"] = 19.17251454462003
$data["G20"] = 3.702166279034726
$data["I20"] = 52.05375058658232
$data["K20"] = 14.17642652515161
$data["B21"] = 11.43468394365751
$data["C21"] = 15.43999553263804
$data["D21"] = 6.579471410145349
$data["E21"] = 20.35469117144591
$data["G21"] = 3.696764115639797
$data["I21"] = 53.95391068655386
$data["K21"] = 14.35667414250225
$data["B22"] = 11.50066307846983
$data["C22"] = 15.81099222526644
$data["D22"] = 6.681235542611192
$data["E22"] = 21.09381076533053
$data["G22"] = 3.693345177149544
$data["I22"] = 55.16120722522145
$data["K22"] = 14.47894801602405
$data["B23"] = 11.46480630468733
$data["C23"] = 15.6138242774589
$data["D23"] = 6.626955941807563
$data["E23"] = 20.70238129746628
$data["G23"] = 3.695159769818562
$data["I23"] = 54.52008418879771
$data["K23"] = 14.41329068662358
$data["B24"] = 11.34319366609603
$data["C24"] = 14.84919152571773
$data["D24"] = 6.420834347067624
$data["E24"] = 19.15262405287433
$data["G24"] = 3.70225647353804
$data["I24"] = 52.02212655850341
$data["K24"] = 14.17355852122423
$data["B25"] = 11.24268951553283
$data["C25"] = 13.9901487939739
$data["D25"] = 6.198630158965843
$data["E25"] = 17.43783091395043
$data["G25"] = 3.710397005347751
$data["I25"] = 49.18728391713557
$data["K25"] = 13.93529155745979

foreach ($cellRef in $data.Keys) {
    $ws.Range($cellRef).Value = $data[$cellRef]
}
